$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price column (D) updates - force text to avoid numeric auto-conversion,
# then restore default style so no stray formatting is introduced.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "34.500.77"
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.794.40"
$ws.Range("D3").Style = "Normal"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "226.81"
$ws.Range("D5").Style = "Normal"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "32.63"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.297"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0695"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0949"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "2.053.31"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "11.05"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.784.14"
$ws.Range("D14").Style = "Normal"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "34.462.77"
$ws.Range("D16").Style = "Normal"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "68.90"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "247.53"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0801"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.23"
$ws.Range("D21").Style = "Normal"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.18"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.07"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "164.21"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.26"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.54"
$ws.Range("D27").Style = "Normal"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.82"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.24"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0522"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.89"
$ws.Range("D33").Style = "Normal"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.437.34"
$ws.Range("D35").Style = "Normal"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "84.62"
$ws.Range("D40").Style = "Normal"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.937"
$ws.Range("D42").Style = "Normal"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0527"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "6.12"
$ws.Range("D46").Style = "Normal"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "105.86"
$ws.Range("D49").Style = "Normal"

# Volume(1h) column (E) updates - plain text assignment (never numeric-looking).
$ws.Range("E2").Value = "  +0.84%  "
$ws.Range("E3").Value = "  +0.20%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("E5").Value = "  +0.05%  "
$ws.Range("E6").Value = "  +1.65%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("E8").Value = "  +2.07%  "
$ws.Range("E9").Value = "  +1.39%  "
$ws.Range("E10").Value = "  +0.77%  "
$ws.Range("E11").Value = "  +0.37%  "
$ws.Range("E12").Value = "  +0.17%  "
$ws.Range("E13").Value = "  -0.77%  "
$ws.Range("E14").Value = "  -0.23%  "
$ws.Range("E15").Value = "  +2.63%  "
$ws.Range("E16").Value = "  +0.89%  "
$ws.Range("E17").Value = "  +2.20%  "
$ws.Range("E18").Value = "  +1.16%  "
$ws.Range("E19").Value = "  +0.79%  "
$ws.Range("E20").Value = "  +2.74%  "
$ws.Range("E21").Value = "  +3.70%  "
$ws.Range("E22").Value = "  +0.10%  "
$ws.Range("E23").Value = "  +1.51%  "
$ws.Range("E24").Value = "  +1.37%  "
$ws.Range("E25").Value = "  +1.97%  "
$ws.Range("E26").Value = "  +1.17%  "
$ws.Range("E27").Value = "  +1.07%  "
$ws.Range("E28").Value = "  +2.23%  "
$ws.Range("E29").Value = "  +0.04%  "
$ws.Range("E30").Value = "  +3.95%  "
$ws.Range("E31").Value = "  +0.09%  "
$ws.Range("E32").Value = "  +0.53%  "
$ws.Range("E33").Value = "  +6.86%  "
$ws.Range("E34").Value = "  +0.95%  "
$ws.Range("E35").Value = "  -1.10%  "
$ws.Range("E36").Value = "  +6.89%  "
$ws.Range("E37").Value = "  +2.99%  "
$ws.Range("E38").Value = "  +3.04%  "
$ws.Range("E39").Value = "  -0.05%  "
$ws.Range("E40").Value = "  +5.20%  "
$ws.Range("E41").Value = "  +1.32%  "
$ws.Range("E42").Value = "  +1.75%  "
$ws.Range("E43").Value = "  +2.06%  "
$ws.Range("E44").Value = "  -0.27%  "
$ws.Range("E45").Value = "  +3.63%  "
$ws.Range("E46").Value = "  +1.15%  "
$ws.Range("E47").Value = "  +0.16%  "
$ws.Range("E48").Value = "  -0.01%  "
$ws.Range("E49").Value = "  -0.08%  "
$ws.Range("E50").Value = "  -3.34%  "
$ws.Range("E51").Value = "  +0.03%  "
